$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Boolean prediction-success flags that flip from TRUE to FALSE
$ws.Range("C8").Value = $false
$ws.Range("C9").Value = $false
$ws.Range("C10").Value = $false
$ws.Range("C18").Value = $false
$ws.Range("C19").Value = $false
$ws.Range("C20").Value = $false

# Updated numeric results (prediction / error / cross-entropy loss / success %)
$ws.Range("D2").Value = 0.9999687317203518
$ws.Range("E2").Value = 0.9999687317203518
$ws.Range("D3").Value = 0.9999985578563054
$ws.Range("E3").Value = 0.9999985578563054
$ws.Range("D4").Value = 0.3441103626345411
$ws.Range("E4").Value = 0.3441103626345411
$ws.Range("D5").Value = 0.9999999999999747
$ws.Range("E5").Value = 0.9999999999999747
$ws.Range("D6").Value = 0.9911316181822624
$ws.Range("E6").Value = 0.9911316181822624
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.000000107814448830965
$ws.Range("E8").Value = 0.9999998921855512
$ws.Range("D9").Value = 0.06323476730530066
$ws.Range("E9").Value = 0.9367652326946994
$ws.Range("D10").Value = 0.00000001020623221607705
$ws.Range("E10").Value = 0.9999999897937678
$ws.Range("D11").Value = 0.9999993037422776
$ws.Range("E11").Value = 0.0000006962577223657007
$ws.Range("F11").Value = 9.747716903686523
$ws.Range("G11").Value = 0.3
$ws.Range("D12").Value = 0.9999986557477916
$ws.Range("E12").Value = 0.9999986557477916
$ws.Range("D13").Value = 0.9999976897315318
$ws.Range("E13").Value = 0.9999976897315318
$ws.Range("D14").Value = 0.02720766984815618
$ws.Range("E14").Value = 0.02720766984815618
$ws.Range("D16").Value = 0.9985188257395833
$ws.Range("E16").Value = 0.9985188257395833
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("D18").Value = 0.000000001075390071267258
$ws.Range("E18").Value = 0.9999999989246099
$ws.Range("D19").Value = 0.04011669241917602
$ws.Range("E19").Value = 0.959883307580824
$ws.Range("D20").Value = 0.0000000201179074517676
$ws.Range("E20").Value = 0.9999999798820925
$ws.Range("D21").Value = 0.9999999818319812
$ws.Range("E21").Value = 0.00000001816801875609286
$ws.Range("F21").Value = 12.58650588989258
$ws.Range("G21").Value = 0.3

Write-Host "Updated classification result values"
